$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "(  11:13:06 AM,  12:24:55 PM, 1:11:49) \n "
$ws.Range("D3").Value = "(  10:57:43 AM,  10:59:07 AM, 0:01:24) \n (  10:59:07 AM,  12:24:29 PM, 1:25:22) \n "
$ws.Range("D4").Value = "(  11:02:50 AM,  11:03:01 AM, 0:00:11) \n "
$ws.Range("D5").Value = "(  10:59:06 AM,  12:03:49 PM, 1:04:43) \n "
$ws.Range("D6").Value = "(  11:16:00 AM,  12:30:00 PM, 1:14:00) \n "
$ws.Range("D7").Value = "(  10:56:12 AM,  10:59:05 AM, 0:02:53) \n (  10:59:06 AM,  12:25:03 PM, 1:25:57) \n "
$ws.Range("D8").Value = "(  11:04:03 AM,  11:38:12 AM, 0:34:09) \n (  11:51:34 AM,  12:09:48 PM, 0:18:14) \n (  12:18:24 PM,  12:25:10 PM, 0:06:46) \n "
$ws.Range("D9").Value = "(  10:56:12 AM,  10:57:27 AM, 0:01:15) \n (  10:57:47 AM,  10:59:05 AM, 0:01:18) \n (  10:59:06 AM,  11:13:19 AM, 0:14:13) \n "
$ws.Range("D10").Value = "(  11:12:16 AM,  12:30:00 PM, 1:17:44) \n "
$ws.Range("D11").Value = "(  11:07:27 AM,  12:24:57 PM, 1:17:30) \n "
$ws.Range("D12").Value = "(  10:57:43 AM,  10:57:57 AM, 0:00:14) \n (  10:58:14 AM,  10:58:18 AM, 0:00:04) \n (  10:59:29 AM,  12:25:08 PM, 1:25:39) \n "
$ws.Range("D13").Value = "(  10:59:39 AM,  12:24:43 PM, 1:25:04) \n "
$ws.Range("D14").Value = "(  10:56:12 AM,  11:05:05 AM, 0:08:53) \n (  11:03:54 AM,  12:24:36 PM, 1:20:42) \n "
$ws.Range("D15").Value = "(  11:18:06 AM,  11:23:05 AM, 0:04:59) \n (  11:29:12 AM,  11:33:57 AM, 0:04:45) \n "
$ws.Range("D16").Value = "(  11:02:04 AM,  12:24:44 PM, 1:22:40) \n "
$ws.Range("D17").Value = "(  11:07:35 AM,  11:20:47 AM, 0:13:12) \n "
$ws.Range("D18").Value = "(  11:00:36 AM,  11:00:46 AM, 0:00:10) \n (  11:01:10 AM,  12:25:07 PM, 1:23:57) \n "
$ws.Range("D19").Value = "(  11:00:30 AM,  11:01:14 AM, 0:00:44) \n (  11:01:14 AM,  12:25:05 PM, 1:23:51) \n "
$ws.Range("D20").Value = "(  11:01:38 AM,  12:25:28 PM, 1:23:50) \n "
$ws.Range("D21").Value = "(  11:06:54 AM,  12:26:56 PM, 1:20:02) \n "
$ws.Range("D22").Value = "(  11:00:34 AM,  11:00:37 AM, 0:00:03) \n "
$ws.Range("D23").Value = "(  11:10:48 AM,  11:14:27 AM, 0:03:39) \n "
$ws.Range("D24").Value = "(  11:08:13 AM,  12:03:57 PM, 0:55:44) \n "
$ws.Range("D25").Value = "(  11:10:25 AM,  12:25:46 PM, 1:15:21) \n "
$ws.Range("D26").Value = "(  11:01:13 AM,  12:24:31 PM, 1:23:18) \n "
$ws.Range("D27").Value = "(  11:37:48 AM,  12:02:10 PM, 0:24:22) \n "
$ws.Range("D28").Value = "(  11:21:57 AM,  12:25:38 PM, 1:03:41) \n "
$ws.Range("D29").Value = "(  11:03:14 AM,  12:24:55 PM, 1:21:41) \n "
$ws.Range("D30").Value = "(  11:02:14 AM,  12:24:51 PM, 1:22:37) \n "
$ws.Range("D31").Value = "(  10:58:20 AM,  11:01:23 AM, 0:03:03) \n "
$ws.Range("D32").Value = "(  11:39:54 AM,  12:02:22 PM, 0:22:28) \n "
$ws.Range("D33").Value = "(  11:11:32 AM,  11:13:41 AM, 0:02:09) \n (  11:13:57 AM,  11:14:05 AM, 0:00:08) \n (  11:18:32 AM,  11:20:49 AM, 0:02:17) \n (  11:22:47 AM,  12:25:32 PM, 1:02:45) \n "
$ws.Range("D34").Value = "(  10:56:12 AM,  10:56:32 AM, 0:00:20) \n "
$ws.Range("D35").Value = "(  11:00:24 AM,  12:11:53 PM, 1:11:29) \n "
$ws.Range("D36").Value = "(  10:57:09 AM,  10:59:33 AM, 0:02:24) \n "
$ws.Range("D37").Value = "(  10:56:12 AM,  10:59:06 AM, 0:02:54) \n (  10:59:06 AM,  12:24:54 PM, 1:25:48) \n "
$ws.Range("D38").Value = "(  10:57:23 AM,  10:59:05 AM, 0:01:42) \n (  10:59:06 AM,  12:25:05 PM, 1:25:59) \n "
$ws.Range("D39").Value = "(  11:50:25 AM,  11:51:45 AM, 0:01:20) \n "
$ws.Range("D40").Value = "(  11:00:06 AM,  11:00:23 AM, 0:00:17) \n "
$ws.Range("D41").Value = "(  11:00:03 AM,  12:19:15 PM, 1:19:12) \n "
$ws.Range("D42").Value = "(  11:00:56 AM,  11:01:00 AM, 0:00:04) \n (  11:01:10 AM,  11:07:54 AM, 0:06:44) \n (  11:10:56 AM,  11:13:12 AM, 0:02:16) \n "
$ws.Range("D43").Value = "(  10:56:12 AM,  12:25:04 PM, 1:28:52) \n (  12:26:13 PM,  12:26:30 PM, 0:00:17) \n "
$ws.Range("D44").Value = "(  10:57:54 AM,  10:58:13 AM, 0:00:19) \n (  10:58:20 AM,  10:59:05 AM, 0:00:45) \n (  10:59:06 AM,  12:25:01 PM, 1:25:55) \n "
$ws.Range("D45").Value = "(  10:56:12 AM,  10:56:18 AM, 0:00:06) \n (  10:56:27 AM,  10:57:22 AM, 0:00:55) \n (  10:58:30 AM,  10:59:04 AM, 0:00:34) \n (  10:59:05 AM,  12:25:02 PM, 1:25:57) \n "
$ws.Range("D46").Value = "(  11:03:18 AM,  11:03:24 AM, 0:00:06) \n (  11:03:24 AM,  11:32:11 AM, 0:28:47) \n (  11:36:23 AM,  12:05:40 PM, 0:29:17) \n (  12:05:51 PM,  12:25:01 PM, 0:19:10) \n "
$ws.Range("D47").Value = "(  11:06:47 AM,  12:25:04 PM, 1:18:17) \n "
$ws.Range("D48").Value = "(  11:04:20 AM,  11:05:16 AM, 0:00:56) \n "
$ws.Range("D49").Value = "(  11:00:26 AM,  11:01:14 AM, 0:00:48) \n (  11:01:15 AM,  12:16:41 PM, 1:15:26) \n "
$ws.Range("D50").Value = "(  11:00:03 AM,  12:30:00 PM, 1:29:57) \n "
$ws.Range("D51").Value = "(  11:01:42 AM,  11:02:18 AM, 0:00:36) \n (  11:03:42 AM,  11:41:47 AM, 0:38:05) \n "
$ws.Range("D52").Value = "(  11:08:36 AM,  11:21:07 AM, 0:12:31) \n "
$ws.Range("D53").Value = "(  11:42:16 AM,  12:25:13 PM, 0:42:57) \n "
$ws.Range("D54").Value = "(  12:13:37 PM,  12:30:00 PM, 0:16:23) \n "
$ws.Range("D55").Value = "(  11:24:22 AM,  11:26:31 AM, 0:02:09) \n "
$ws.Range("D56").Value = "(  11:16:24 AM,  12:30:00 PM, 1:13:36) \n "
$ws.Range("D57").Value = "(  11:01:15 AM,  11:01:32 AM, 0:00:17) \n (  11:01:32 AM,  12:24:32 PM, 1:23:00) \n "
$ws.Range("D58").Value = "(  11:59:35 AM,  12:02:09 PM, 0:02:34) \n (  12:02:53 PM,  12:30:00 PM, 0:27:07) \n "
$ws.Range("D59").Value = "(  10:59:29 AM,  11:02:18 AM, 0:02:49) \n (  11:17:23 AM,  12:25:50 PM, 1:08:27) \n "
$ws.Range("D60").Value = "(  10:56:32 AM,  10:59:04 AM, 0:02:32) \n (  10:59:05 AM,  12:25:13 PM, 1:26:08) \n "
$ws.Range("D61").Value = "(  11:07:05 AM,  12:01:45 PM, 0:54:40) \n "
$ws.Range("D62").Value = "(  10:58:54 AM,  10:59:04 AM, 0:00:10) \n (  10:59:05 AM,  12:18:38 PM, 1:19:33) \n (  12:22:20 PM,  12:30:00 PM, 0:07:40) \n "
$ws.Range("D63").Value = "(  10:57:44 AM,  10:59:04 AM, 0:01:20) \n (  10:59:06 AM,  11:36:12 AM, 0:37:06) \n (  11:37:03 AM,  12:26:07 PM, 0:49:04) \n "
$ws.Range("D64").Value = "(  11:12:09 AM,  11:36:04 AM, 0:23:55) \n "
$ws.Range("D65").Value = "(  11:01:23 AM,  12:30:00 PM, 1:28:37) \n "
$ws.Range("D66").Value = "(  10:57:51 AM,  10:58:22 AM, 0:00:31) \n (  10:58:57 AM,  10:59:03 AM, 0:00:06) \n (  10:59:05 AM,  12:25:48 PM, 1:26:43) \n "
$ws.Range("D67").Value = "(  11:01:42 AM,  11:01:49 AM, 0:00:07) \n (  11:03:06 AM,  11:03:13 AM, 0:00:07) \n (  11:04:26 AM,  11:31:12 AM, 0:26:46) \n (  11:32:21 AM,  12:21:46 PM, 0:49:25) \n "
$ws.Range("D68").Value = "(  11:08:35 AM,  12:24:27 PM, 1:15:52) \n "
$ws.Range("D69").Value = "(  11:00:01 AM,  11:00:06 AM, 0:00:05) \n (  11:00:22 AM,  12:25:36 PM, 1:25:14) \n "
$ws.Range("D70").Value = "(  11:50:48 AM,  12:13:12 PM, 0:22:24) \n "
$ws.Range("D71").Value = "(  10:56:12 AM,  10:56:16 AM, 0:00:04) \n (  11:01:29 AM,  11:01:37 AM, 0:00:08) \n (  11:03:05 AM,  11:03:11 AM, 0:00:06) \n (  11:15:20 AM,  12:16:00 PM, 1:00:40) \n "
$ws.Range("D72").Value = "(  11:00:09 AM,  11:02:57 AM, 0:02:48) \n (  11:03:32 AM,  12:25:07 PM, 1:21:35) \n "
$ws.Range("D73").Value = "(  11:51:10 AM,  12:24:55 PM, 0:33:45) \n "
$ws.Range("D74").Value = "(  11:11:17 AM,  12:04:07 PM, 0:52:50) \n "
$ws.Range("D75").Value = "(  11:07:51 AM,  12:25:13 PM, 1:17:22) \n "
$ws.Range("D76").Value = "(  11:00:14 AM,  12:30:00 PM, 1:29:46) \n "
$ws.Range("D77").Value = "(  11:16:55 AM,  11:48:32 AM, 0:31:37) \n (  11:48:46 AM,  11:49:33 AM, 0:00:47) \n "
$ws.Range("D78").Value = "(  10:56:12 AM,  10:56:34 AM, 0:00:22) \n "
$ws.Range("D79").Value = "(  11:02:03 AM,  12:30:00 PM, 1:27:57) \n "
$ws.Range("D80").Value = "(  11:04:31 AM,  11:04:45 AM, 0:00:14) \n "
$ws.Range("D81").Value = "(  10:57:08 AM,  10:59:05 AM, 0:01:57) \n (  10:59:05 AM,  12:21:52 PM, 1:22:47) \n "
$ws.Range("D82").Value = "(  11:05:59 AM,  11:06:25 AM, 0:00:26) \n (  11:06:40 AM,  11:51:48 AM, 0:45:08) \n (  11:57:55 AM,  12:30:00 PM, 0:32:05) \n "
$ws.Range("D83").Value = "(  12:21:14 PM,  12:25:10 PM, 0:03:56) \n "
$ws.Range("D84").Value = "(  10:57:44 AM,  10:59:04 AM, 0:01:20) \n (  10:59:04 AM,  11:55:52 AM, 0:56:48) \n (  11:59:19 AM,  12:24:53 PM, 0:25:34) \n "
$ws.Range("D85").Value = "(  11:32:42 AM,  12:24:26 PM, 0:51:44) \n "
$ws.Range("D86").Value = "(  11:06:32 AM,  11:23:56 AM, 0:17:24) \n "
$ws.Range("D87").Value = "(  11:02:11 AM,  12:30:00 PM, 1:27:49) \n "
$ws.Range("D88").Value = "(  11:01:49 AM,  11:01:58 AM, 0:00:09) \n (  11:01:58 AM,  11:03:24 AM, 0:01:26) \n (  11:03:28 AM,  12:25:03 PM, 1:21:35) \n "
$ws.Range("D89").Value = "(  11:31:47 AM,  11:41:58 AM, 0:10:11) \n (  11:59:31 AM,  12:04:22 PM, 0:04:51) \n "
$ws.Range("D90").Value = "(  10:59:31 AM,  12:25:00 PM, 1:25:29) \n "
$ws.Range("D91").Value = "(  11:00:46 AM,  12:24:38 PM, 1:23:52) \n "
$ws.Range("D92").Value = "(  11:11:15 AM,  12:30:00 PM, 1:18:45) \n "
$ws.Range("D93").Value = "(  11:01:01 AM,  11:01:13 AM, 0:00:12) \n (  11:01:13 AM,  12:24:29 PM, 1:23:16) \n "
$ws.Range("D94").Value = "(  11:16:33 AM,  11:18:41 AM, 0:02:08) \n "
$ws.Range("D95").Value = "(  11:00:22 AM,  12:24:41 PM, 1:24:19) \n "
$ws.Range("D96").Value = "(  10:59:19 AM,  12:24:23 PM, 1:25:04) \n "
$ws.Range("D97").Value = "(  11:00:31 AM,  11:47:23 AM, 0:46:52) \n "
$ws.Range("D98").Value = "(  10:56:12 AM,  10:59:03 AM, 0:02:51) \n (  10:57:29 AM,  10:59:05 AM, 0:01:36) \n (  10:59:04 AM,  11:37:00 AM, 0:37:56) \n (  10:59:06 AM,  12:25:18 PM, 1:26:12) \n "
$ws.Range("D99").Value = "(  10:56:12 AM,  12:25:15 PM, 1:29:03) \n "
$ws.Range("D100").Value = "(  11:06:11 AM,  12:05:09 PM, 0:58:58) \n "
$ws.Range("D101").Value = "(  11:01:34 AM,  11:12:49 AM, 0:11:15) \n "
$ws.Range("D102").Value = "(  11:45:23 AM,  12:25:08 PM, 0:39:45) \n "
$ws.Range("D103").Value = "(  11:00:09 AM,  12:25:02 PM, 1:24:53) \n "
$ws.Range("D104").Value = "(  10:56:14 AM,  10:56:21 AM, 0:00:07) \n (  10:59:45 AM,  11:00:21 AM, 0:00:36) \n "
$ws.Range("D105").Value = "(  11:00:37 AM,  11:01:15 AM, 0:00:38) \n (  11:01:15 AM,  11:23:00 AM, 0:21:45) \n (  12:12:13 PM,  12:22:38 PM, 0:10:25) \n "
$ws.Range("D106").Value = "(  11:02:06 AM,  12:25:00 PM, 1:22:54) \n "
$ws.Range("D107").Value = "(  11:19:57 AM,  11:21:25 AM, 0:01:28) \n (  11:24:00 AM,  12:25:03 PM, 1:01:03) \n "
$ws.Range("D108").Value = "(  10:56:14 AM,  10:59:06 AM, 0:02:52) \n (  10:59:07 AM,  11:59:05 AM, 0:59:58) \n "
$ws.Range("D109").Value = "(  10:59:53 AM,  12:25:04 PM, 1:25:11) \n "
$ws.Range("D110").Value = "(  11:02:49 AM,  11:03:24 AM, 0:00:35) \n (  11:03:24 AM,  12:24:58 PM, 1:21:34) \n "
$ws.Range("D111").Value = "(  10:57:56 AM,  10:59:07 AM, 0:01:11) \n (  10:59:07 AM,  12:24:36 PM, 1:25:29) \n (  11:34:46 AM,  12:27:15 PM, 0:52:29) \n "
$ws.Range("D112").Value = "(  10:57:39 AM,  10:59:07 AM, 0:01:28) \n (  10:59:07 AM,  12:24:36 PM, 1:25:29) \n "
$ws.Range("D113").Value = "(  11:21:29 AM,  12:24:58 PM, 1:03:29) \n "
$ws.Range("D114").Value = "(  10:59:25 AM,  12:17:44 PM, 1:18:19) \n "
$ws.Range("D115").Value = "(  11:02:34 AM,  11:34:36 AM, 0:32:02) \n (  11:35:02 AM,  12:25:04 PM, 0:50:02) \n "
$ws.Range("D116").Value = "(  11:26:42 AM,  12:25:05 PM, 0:58:23) \n "
$ws.Range("D117").Value = "(  11:00:44 AM,  11:01:14 AM, 0:00:30) \n (  11:01:14 AM,  11:21:55 AM, 0:20:41) \n "
$ws.Range("D118").Value = "(  11:04:26 AM,  11:04:35 AM, 0:00:09) \n "
$ws.Range("D119").Value = "(  11:02:49 AM,  12:24:58 PM, 1:22:09) \n "
$ws.Range("D120").Value = "(  11:05:04 AM,  11:05:46 AM, 0:00:42) \n (  11:05:46 AM,  12:30:00 PM, 1:24:14) \n "
$ws.Range("D121").Value = "(  11:20:41 AM,  11:23:17 AM, 0:02:36) \n "
$ws.Range("D122").Value = "(  10:59:13 AM,  10:59:21 AM, 0:00:08) \n "
$ws.Range("D123").Value = "(  11:24:21 AM,  12:24:59 PM, 1:00:38) \n "
$ws.Range("D124").Value = "(  10:58:03 AM,  10:58:08 AM, 0:00:05) \n (  10:59:07 AM,  12:24:52 PM, 1:25:45) \n "
$ws.Range("D125").Value = "(  11:32:38 AM,  12:25:02 PM, 0:52:24) \n "
$ws.Range("D126").Value = "(  10:56:12 AM,  10:57:10 AM, 0:00:58) \n (  11:00:13 AM,  12:14:51 PM, 1:14:38) \n (  12:16:25 PM,  12:30:00 PM, 0:13:35) \n "
$ws.Range("D127").Value = "(  11:08:45 AM,  11:14:40 AM, 0:05:55) \n "
$ws.Range("D128").Value = "(  10:58:43 AM,  10:59:04 AM, 0:00:21) \n (  10:59:05 AM,  12:25:10 PM, 1:26:05) \n "
$ws.Range("D129").Value = "(  11:35:43 AM,  12:23:35 PM, 0:47:52) \n "
$ws.Range("D130").Value = "(  11:13:09 AM,  11:50:42 AM, 0:37:33) \n "
$ws.Range("D131").Value = "(  10:57:37 AM,  10:59:03 AM, 0:01:26) \n (  10:59:04 AM,  12:30:00 PM, 1:30:56) \n "
$ws.Range("D132").Value = "(  11:02:33 AM,  12:24:56 PM, 1:22:23) \n "
$ws.Range("D133").Value = "(  11:00:51 AM,  12:25:04 PM, 1:24:13) \n "
$ws.Range("D134").Value = "(  11:06:52 AM,  12:25:05 PM, 1:18:13) \n "
$ws.Range("D135").Value = "(  10:57:47 AM,  10:58:09 AM, 0:00:22) \n "
$ws.Range("D136").Value = "(  11:01:49 AM,  11:37:50 AM, 0:36:01) \n (  11:38:00 AM,  12:16:12 PM, 0:38:12) \n (  12:17:57 PM,  12:25:14 PM, 0:07:17) \n "
$ws.Range("D137").Value = "(  11:09:46 AM,  11:46:00 AM, 0:36:14) \n (  11:46:14 AM,  12:25:16 PM, 0:39:02) \n "
$ws.Range("D138").Value = "(  11:05:06 AM,  12:24:57 PM, 1:19:51) \n "
$ws.Range("D139").Value = "(  10:58:49 AM,  10:59:05 AM, 0:00:16) \n "
$ws.Range("D140").Value = "(  11:13:00 AM,  11:15:15 AM, 0:02:15) \n "
$ws.Range("D141").Value = "(  12:22:04 PM,  12:24:57 PM, 0:02:53) \n "
$ws.Range("D142").Value = "(  11:01:57 AM,  11:37:08 AM, 0:35:11) \n "
$ws.Range("D143").Value = "(  11:08:28 AM,  12:24:45 PM, 1:16:17) \n "
$ws.Range("D144").Value = "(  11:07:24 AM,  11:56:41 AM, 0:49:17) \n "
$ws.Range("D145").Value = "(  11:02:14 AM,  12:25:03 PM, 1:22:49) \n "
$ws.Range("D146").Value = "(  11:00:01 AM,  11:05:38 AM, 0:05:37) \n "
$ws.Range("D147").Value = "(  10:56:12 AM,  12:25:05 PM, 1:28:53) \n "
$ws.Range("D148").Value = "(  11:00:27 AM,  11:01:14 AM, 0:00:47) \n (  11:01:14 AM,  12:25:06 PM, 1:23:52) \n "
$ws.Range("D149").Value = "(  11:07:16 AM,  11:15:45 AM, 0:08:29) \n (  11:16:35 AM,  12:25:01 PM, 1:08:26) \n "
$ws.Range("D150").Value = "(  11:40:22 AM,  11:41:35 AM, 0:01:13) \n "
$ws.Range("D151").Value = "(  11:01:33 AM,  12:03:59 PM, 1:02:26) \n "
$ws.Range("D152").Value = "(  11:00:08 AM,  12:27:30 PM, 1:27:22) \n "
$ws.Range("D153").Value = "(  11:01:47 AM,  12:24:35 PM, 1:22:48) \n "
$ws.Range("D154").Value = "(  10:56:53 AM,  10:57:01 AM, 0:00:08) \n (  10:57:07 AM,  10:59:05 AM, 0:01:58) \n (  10:58:40 AM,  10:59:05 AM, 0:00:25) \n (  10:59:06 AM,  11:00:25 AM, 0:01:19) \n (  11:01:08 AM,  12:26:59 PM, 1:25:51) \n "
$ws.Range("D155").Value = "(  10:58:27 AM,  10:58:52 AM, 0:00:25) \n (  11:00:32 AM,  12:30:00 PM, 1:29:28) \n "
$ws.Range("D156").Value = "(  11:15:27 AM,  12:25:13 PM, 1:09:46) \n "
$ws.Range("D157").Value = "(  10:56:16 AM,  10:59:05 AM, 0:02:49) \n (  10:59:06 AM,  12:24:58 PM, 1:25:52) \n "
$ws.Range("D158").Value = "(  11:03:01 AM,  12:26:57 PM, 1:23:56) \n "
$ws.Range("D159").Value = "(  10:58:58 AM,  10:59:33 AM, 0:00:35) \n (  11:00:30 AM,  12:03:51 PM, 1:03:21) \n "
$ws.Range("D160").Value = "(  11:14:25 AM,  11:20:25 AM, 0:06:00) \n "
$ws.Range("D161").Value = "(  11:03:40 AM,  12:14:23 PM, 1:10:43) \n (  12:14:30 PM,  12:25:04 PM, 0:10:34) \n "
$ws.Range("D162").Value = "(  10:59:47 AM,  11:00:02 AM, 0:00:15) \n (  11:00:02 AM,  12:25:12 PM, 1:25:10) \n "
$ws.Range("D163").Value = "(  11:11:30 AM,  12:09:29 PM, 0:57:59) \n "
$ws.Range("D164").Value = "(  11:07:40 AM,  12:30:00 PM, 1:22:20) \n "
$ws.Range("D165").Value = "(  11:00:05 AM,  12:15:24 PM, 1:15:19) \n "
$ws.Range("D166").Value = "(  10:58:08 AM,  10:59:05 AM, 0:00:57) \n (  10:59:06 AM,  12:24:29 PM, 1:25:23) \n "
$ws.Range("D167").Value = "(  11:10:35 AM,  12:25:16 PM, 1:14:41) \n "
$ws.Range("D168").Value = "(  11:04:00 AM,  11:04:13 AM, 0:00:13) \n "
$ws.Range("D169").Value = "(  10:59:12 AM,  11:16:55 AM, 0:17:43) \n (  11:17:42 AM,  12:24:01 PM, 1:06:19) \n "
$ws.Range("D170").Value = "(  11:27:38 AM,  11:28:47 AM, 0:01:09) \n (  11:30:34 AM,  11:44:52 AM, 0:14:18) \n "
$ws.Range("D171").Value = "(  11:02:43 AM,  11:02:59 AM, 0:00:16) \n (  11:03:19 AM,  12:25:01 PM, 1:21:42) \n "
$ws.Range("D172").Value = "(  11:55:42 AM,  11:58:13 AM, 0:02:31) \n "
$ws.Range("D173").Value = "(  11:04:42 AM,  11:04:59 AM, 0:00:17) \n (  11:04:59 AM,  11:05:34 AM, 0:00:35) \n (  11:07:32 AM,  11:12:04 AM, 0:04:32) \n "
$ws.Range("D174").Value = "(  11:05:07 AM,  11:14:35 AM, 0:09:28) \n "
$ws.Range("D175").Value = "(  11:16:44 AM,  11:45:33 AM, 0:28:49) \n "
$ws.Range("D176").Value = "(  10:56:12 AM,  12:25:10 PM, 1:28:58) \n "
$ws.Range("D177").Value = "(  10:56:19 AM,  12:25:34 PM, 1:29:15) \n "
$ws.Range("D178").Value = "(  11:03:54 AM,  12:23:36 PM, 1:19:42) \n "
$ws.Range("D179").Value = "(  11:56:03 AM,  12:06:08 PM, 0:10:05) \n "
$ws.Range("D180").Value = "(  11:00:47 AM,  11:01:11 AM, 0:00:24) \n "
$ws.Range("D181").Value = "(  11:12:32 AM,  12:03:03 PM, 0:50:31) \n (  12:04:41 PM,  12:11:01 PM, 0:06:20) \n "
$ws.Range("D182").Value = "(  11:01:03 AM,  11:05:09 AM, 0:04:06) \n (  11:05:12 AM,  11:07:02 AM, 0:01:50) \n "
$ws.Range("D183").Value = "(  11:45:53 AM,  11:48:00 AM, 0:02:07) \n "
$ws.Range("D184").Value = "(  11:01:44 AM,  12:16:50 PM, 1:15:06) \n "
$ws.Range("D185").Value = "(  11:07:04 AM,  12:24:40 PM, 1:17:36) \n "
$ws.Range("D186").Value = "(  11:01:17 AM,  11:03:05 AM, 0:01:48) \n (  11:02:50 AM,  11:03:22 AM, 0:00:32) \n (  11:03:22 AM,  11:07:12 AM, 0:03:50) \n (  11:06:27 AM,  12:15:03 PM, 1:08:36) \n (  12:13:07 PM,  12:25:08 PM, 0:12:01) \n "
